# Weekly fruit/vegetable data update:
# Insert a new weekly observation row at row 29 (pushing the existing
# rows 29-129 down to 30-130), and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("29:29").Insert()

$ws.Range("A29").Value = 9
$ws.Range("B29").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44525
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100101
$ws.Range("H29").Value = "Berries"
$ws.Range("I29").Value = 100101001
$ws.Range("J29").Value = "Arándano (blue)"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 500
$ws.Range("N29").Value = 5000
$ws.Range("O29").Value = 5000
$ws.Range("P29").Value = 5000
$ws.Range("Q29").Value = "$/bandeja 2 kilos"
$ws.Range("R29").Value = "Provincia de Curicó"
$ws.Range("S29").Value = 2500
$ws.Range("T29").Value = 2
